# Update Sage scrape results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Offensive Cyber Operations and State Power: Lessons from Russia in Ukraine"
$ws.Range("E2").Value = "10.1177/00207020241234228"

# Row 3
$ws.Range("B3").Value = "Cyber scares and prophylactic policies: Crossnational evidence on the effect of cyberattacks on public support for surveillance"
$ws.Range("E3").Value = "10.1177/00223433241233960"
$ws.Range("F3").Value = "Restricted"

# Row 4
$ws.Range("B4").Value = "Cyberattacks and public opinion – The effect of uncertainty in guiding preferences"
$ws.Range("E4").Value = "10.1177/00223433231218178"

# Row 5
$ws.Range("B5").Value = "Moving beyond the sanctuary paradigm: Canada must face up to the reality of a contested and dangerous space environment"
$ws.Range("E5").Value = "10.1177/00207020231178394"

# Row 6
$ws.Range("B6").Value = "Cyber-Flirting: Playing at Love on the Internet"
$ws.Range("E6").Value = "10.1177/0959354303013003003"

# Row 7
$ws.Range("B7").Value = "Considering the cost of cyber warfare: advancing cyber warfare analytics to better assess tradeoffs in system destruction warfare"
$ws.Range("E7").Value = "10.1177/15485129221114354"

# Row 8
$ws.Range("B8").Value = "ROBIN: An open-source middleware for plug‘n’produce of Cyber-Physical Systems"
$ws.Range("E8").Value = "10.1177/1729881420910316"
$ws.Range("F8").Value = "Open Access"

# Row 9
$ws.Range("B9").Value = "Securing Virtual Space: Cyber War, Cyber Terror, and Risk"
$ws.Range("E9").Value = "10.1177/1206331211430016"
$ws.Range("F9").Value = "Restricted"

# Row 10
$ws.Range("B10").Value = "Tech titans, cyber commons and the war in Ukraine: An incipient shift in international relations"
$ws.Range("E10").Value = "10.1177/00471178231211500"

# Row 11
$ws.Range("B11").Value = "Invisible Digital Front: Can Cyber Attacks Shape Battlefield Events?"
$ws.Range("E11").Value = "10.1177/0022002717737138"
